# Fill in the bill form for "Mr. Sunanda Das" (label-wise individual bill)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header block: Name / Designation / Department
$ws.Range("A3").Value = "নাম: Mr. Sunanda Das"
$ws.Range("A4").Value = "পদবী: সহকারী অধ্যাপক"
$ws.Range("F5").Value = "বিভাগ :সিএসই"

# Quantity entries for the bill lines
$ws.Range("G16").Value = 27
$ws.Range("G20").Value = 40
$ws.Range("G26").Value = 1
$ws.Range("G29").Value = 1

# Amount in words
$ws.Range("A32").Value = "কথায়:সাত হাজার দু'শ সত্তর সাতটি টাকা মাত্র।"

# Move the active selection to B5 (as last saved in the workbook)
$ws.Range("B5").Select() | Out-Null
